# Edit the "Rectangle 5" shape (Azure PaaS noSQL Database box) on slide 1.
# Its text is:
#   Azure PaaS noSQL Database
#   (Table storage, DocumentDB, HDInsight Hbase, …)
# and must become:
#   Azure PaaS noSQL Database
#   (DocumentDB, HDInsight Hbase, Table storage, …)
#
# The second line is made of five runs (after the line break):
#   R1 "(Table "                    rPr: lang="fr-FR" dirty="0"
#   R2 "storage"                    rPr: lang="fr-FR" dirty="0" err="1"
#   R3 ", DocumentDB, HDInsight "   rPr: lang="fr-FR" dirty="0"
#   R4 "Hbase"                      rPr: lang="fr-FR" dirty="0" err="1"
#   R5 ", …)"                       rPr: lang="fr-FR"
#
# and must become seven runs:
#   "(DocumentDB, HDInsight "  (R1's formatting)
#   "Hbase"                    (R2's formatting, err="1")
#   ", "                       (R3's formatting)
#   "Table "                  (R3's formatting)
#   "storage"                  (R4's formatting, err="1")
#   ", "                       (R5's formatting)
#   "…)"                       (R5's formatting)
#
# Rather than physically moving runs around (which is not exposed by this
# object model), each original run is edited *in place*, keeping its own
# a:rPr untouched, so the relative left-to-right run order is preserved
# while the rendered text ends up rearranged as required. Edits are applied
# right-to-left (highest character offset first) so earlier replacements
# never invalidate the character offsets used by replacements still to come.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Rectangle 5") {
        $shape = $candidate
    }
}
if ($shape -eq $null) {
    # Fallback: this shape is the third one added to the slide.
    $shape = $s.Shapes.Item(3)
}

$tr = $shape.TextFrame.TextRange

# Locate the second paragraph line (after the manual line break) by finding
# the opening parenthesis that starts "(Table storage, ...".
$full = $tr.Text
$parenPos = $full.IndexOf("(Table")
$lineStart = $parenPos + 1   # 1-based character index of "("

# Original run boundaries, relative to $lineStart:
#   R1 "(Table "                  -> offset 0,  length 7
#   R2 "storage"                  -> offset 7,  length 7
#   R3 ", DocumentDB, HDInsight " -> offset 14, length 24
#   R4 "Hbase"                    -> offset 38, length 5
#   R5 ", \u2026)"                -> offset 43, length 4
$r1Start = $lineStart
$r2Start = $lineStart + 7
$r3Start = $lineStart + 14
$r4Start = $lineStart + 38
$r5Start = $lineStart + 43

$ellipsis = [char]0x2026

# --- R5 ", …)" -> split in place into ", " + "…)" (same total length: 4) ---
$r5a = $tr.Characters($r5Start, 2)
$r5a.Text = ", "
$r5b = $tr.Characters($r5Start + 2, 2)
$r5b.Text = $ellipsis + ")"

# --- R4 "Hbase" -> "storage" (keeps R4's err="1" formatting) ---
$r4 = $tr.Characters($r4Start, 5)
$r4.Text = "storage"

# --- R3 ", DocumentDB, HDInsight " -> ", Table " then split into ", " + "Table " ---
$r3 = $tr.Characters($r3Start, 24)
$r3.Text = ", Table "
$r3a = $tr.Characters($r3Start, 2)
$r3a.Text = ", "
$r3b = $tr.Characters($r3Start + 2, 6)
$r3b.Text = "Table "

# --- R2 "storage" -> "Hbase" (keeps R2's err="1" formatting) ---
$r2 = $tr.Characters($r2Start, 7)
$r2.Text = "Hbase"

# --- R1 "(Table " -> "(DocumentDB, HDInsight " ---
$r1 = $tr.Characters($r1Start, 7)
$r1.Text = "(DocumentDB, HDInsight "

Write-Host "Updated text:" $tr.Text
